$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 4, shifting the existing rows 4:37 down to 5:38.
# Excel copies the formatting of the row above (row 3) onto the new row,
# which preserves the date style (s="2") already present on column D.
$ws.Rows.Item(4).Insert()

# Populate the newly inserted row 4 with the new weekly record.
$ws.Cells.Item(4, 1).Value = 10
$ws.Cells.Item(4, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(4, 3).Value = "La Araucanía"
$ws.Cells.Item(4, 4).Value = [DateTime]"2022-05-18"
$ws.Cells.Item(4, 5).Value = 9
$ws.Cells.Item(4, 6).Value = "Fruta"
$ws.Cells.Item(4, 7).Value = 100108
$ws.Cells.Item(4, 8).Value = "Tropicales y subtropicales"
$ws.Cells.Item(4, 9).Value = 100108003
$ws.Cells.Item(4, 10).Value = "Maracuyá"
$ws.Cells.Item(4, 11).Value = "Sin especificar"
$ws.Cells.Item(4, 12).Value = "Primera"
$ws.Cells.Item(4, 13).Value = 20
$ws.Cells.Item(4, 14).Value = 35000
$ws.Cells.Item(4, 15).Value = 35000
$ws.Cells.Item(4, 16).Value = 35000
$ws.Cells.Item(4, 17).Value = "$/caja 18 kilos"
$ws.Cells.Item(4, 18).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(4, 19).Value = 1944
$ws.Cells.Item(4, 20).Value = 18
